# Apply KHL stats update: add 6 new matches to Matches_SOG, and refresh
# rolling aggregate stats (Shots_HA, Shots_Summary, Meta_ext) to reflect
# games through 2025-12-20T17:00:00Z (commit: chore(runtime): publish files + archive).

$wb = $excel.ActiveWorkbook
$wsMatches = $wb.Worksheets.Item("Matches_SOG")
$wsHA = $wb.Worksheets.Item("Shots_HA")
$wsSS = $wb.Worksheets.Item("Shots_Summary")
$wsMeta = $wb.Worksheets.Item("Meta_ext")

# ---- Matches_SOG: append 6 new match rows (389-394) ----
$wsMatches.Range("A389").Value = "'897881"
$wsMatches.Range("B389").Value = '2025-12-19T15:30:00'
$wsMatches.Range("C389").Value = 'Сибирь'
$wsMatches.Range("D389").Value = 'Адмирал'
$wsMatches.Range("E389").Value = 26
$wsMatches.Range("F389").Value = 36
$wsMatches.Range("G389").Value = 'khl_text'

$wsMatches.Range("A390").Value = "'897883"
$wsMatches.Range("B390").Value = '2025-12-19T17:00:00'
$wsMatches.Range("C390").Value = 'Трактор'
$wsMatches.Range("D390").Value = 'Металлург Мг'
$wsMatches.Range("E390").Value = 38
$wsMatches.Range("F390").Value = 24
$wsMatches.Range("G390").Value = 'khl_text'

$wsMatches.Range("A391").Value = "'897882"
$wsMatches.Range("B391").Value = '2025-12-19T19:30:00'
$wsMatches.Range("C391").Value = 'Драконы'
$wsMatches.Range("D391").Value = 'Лада'
$wsMatches.Range("E391").Value = 46
$wsMatches.Range("F391").Value = 27
$wsMatches.Range("G391").Value = 'khl_text'

$wsMatches.Range("A392").Value = "'897884"
$wsMatches.Range("B392").Value = '2025-12-20T14:30:00'
$wsMatches.Range("C392").Value = 'Салават Юлаев'
$wsMatches.Range("D392").Value = 'Ак Барс'
$wsMatches.Range("E392").Value = 35
$wsMatches.Range("F392").Value = 26
$wsMatches.Range("G392").Value = 'khl_text'

$wsMatches.Range("A393").Value = "'897885"
$wsMatches.Range("B393").Value = '2025-12-20T17:00:00'
$wsMatches.Range("C393").Value = 'СКА'
$wsMatches.Range("D393").Value = 'Спартак'
$wsMatches.Range("E393").Value = 22
$wsMatches.Range("F393").Value = 41
$wsMatches.Range("G393").Value = 'khl_text'

$wsMatches.Range("A394").Value = "'897886"
$wsMatches.Range("B394").Value = '2025-12-20T17:00:00'
$wsMatches.Range("C394").Value = 'Локомотив'
$wsMatches.Range("D394").Value = 'Авангард'
$wsMatches.Range("E394").Value = 33
$wsMatches.Range("F394").Value = 34
$wsMatches.Range("G394").Value = 'khl_text'

# ---- Shots_HA ----
$wsHA.Range("D2").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("F2").Value = 18
$wsHA.Range("K2").Value = 635
$wsHA.Range("L2").Value = 593
$wsHA.Range("M2").Value = 35.3
$wsHA.Range("D3").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("D4").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("F4").Value = 18
$wsHA.Range("K4").Value = 563
$wsHA.Range("L4").Value = 506
$wsHA.Range("M4").Value = 31.3
$wsHA.Range("N4").Value = 28.1
$wsHA.Range("D5").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("F5").Value = 19
$wsHA.Range("K5").Value = 620
$wsHA.Range("L5").Value = 558
$wsHA.Range("M5").Value = 32.6
$wsHA.Range("N5").Value = 29.4
$wsHA.Range("D6").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("D7").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("D8").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("D9").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("D10").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("E10").Value = 15
$wsHA.Range("G10").Value = 454
$wsHA.Range("H10").Value = 515
$wsHA.Range("I10").Value = 30.3
$wsHA.Range("J10").Value = 34.3
$wsHA.Range("D11").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("F11").Value = 19
$wsHA.Range("K11").Value = 462
$wsHA.Range("L11").Value = 729
$wsHA.Range("M11").Value = 24.3
$wsHA.Range("N11").Value = 38.4
$wsHA.Range("D12").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("E12").Value = 20
$wsHA.Range("G12").Value = 659
$wsHA.Range("H12").Value = 530
$wsHA.Range("I12").Value = 33
$wsHA.Range("J12").Value = 26.5
$wsHA.Range("D13").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("F13").Value = 18
$wsHA.Range("K13").Value = 542
$wsHA.Range("L13").Value = 536
$wsHA.Range("M13").Value = 30.1
$wsHA.Range("N13").Value = 29.8
$wsHA.Range("D14").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("D15").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("E15").Value = 19
$wsHA.Range("G15").Value = 605
$wsHA.Range("H15").Value = 657
$wsHA.Range("I15").Value = 31.8
$wsHA.Range("J15").Value = 34.6
$wsHA.Range("D16").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("E16").Value = 15
$wsHA.Range("G16").Value = 424
$wsHA.Range("H16").Value = 415
$wsHA.Range("I16").Value = 28.3
$wsHA.Range("J16").Value = 27.7
$wsHA.Range("D17").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("D18").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("E18").Value = 17
$wsHA.Range("G18").Value = 449
$wsHA.Range("H18").Value = 581
$wsHA.Range("J18").Value = 34.2
$wsHA.Range("D19").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("F19").Value = 17
$wsHA.Range("K19").Value = 549
$wsHA.Range("L19").Value = 548
$wsHA.Range("M19").Value = 32.3
$wsHA.Range("N19").Value = 32.2
$wsHA.Range("D20").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("D21").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("E21").Value = 19
$wsHA.Range("G21").Value = 689
$wsHA.Range("H21").Value = 553
$wsHA.Range("I21").Value = 36.3
$wsHA.Range("J21").Value = 29.1
$wsHA.Range("D22").Value = '2025-12-20T17:00:00Z'
$wsHA.Range("D23").Value = '2025-12-20T17:00:00Z'

# ---- Shots_Summary ----
$wsSS.Range("D2").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("E2").Value = 34
$wsSS.Range("F2").Value = 1149
$wsSS.Range("G2").Value = 1054
$wsSS.Range("I2").Value = 31
$wsSS.Range("D3").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("D4").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("E4").Value = 34
$wsSS.Range("F4").Value = 1147
$wsSS.Range("G4").Value = 939
$wsSS.Range("I4").Value = 27.6
$wsSS.Range("D5").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("E5").Value = 39
$wsSS.Range("F5").Value = 1280
$wsSS.Range("G5").Value = 1112
$wsSS.Range("H5").Value = 32.8
$wsSS.Range("I5").Value = 28.5
$wsSS.Range("D6").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("D7").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("D8").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("D9").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("D10").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("E10").Value = 35
$wsSS.Range("F10").Value = 1020
$wsSS.Range("G10").Value = 1240
$wsSS.Range("H10").Value = 29.1
$wsSS.Range("I10").Value = 35.4
$wsSS.Range("D11").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("E11").Value = 35
$wsSS.Range("F11").Value = 884
$wsSS.Range("G11").Value = 1303
$wsSS.Range("H11").Value = 25.3
$wsSS.Range("I11").Value = 37.2
$wsSS.Range("D12").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("E12").Value = 39
$wsSS.Range("F12").Value = 1235
$wsSS.Range("G12").Value = 1000
$wsSS.Range("H12").Value = 31.7
$wsSS.Range("I12").Value = 25.6
$wsSS.Range("D13").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("E13").Value = 35
$wsSS.Range("F13").Value = 1184
$wsSS.Range("G13").Value = 966
$wsSS.Range("H13").Value = 33.8
$wsSS.Range("I13").Value = 27.6
$wsSS.Range("D14").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("D15").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("E15").Value = 36
$wsSS.Range("F15").Value = 1118
$wsSS.Range("G15").Value = 1240
$wsSS.Range("H15").Value = 31.1
$wsSS.Range("I15").Value = 34.4
$wsSS.Range("D16").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("E16").Value = 37
$wsSS.Range("F16").Value = 1055
$wsSS.Range("G16").Value = 1097
$wsSS.Range("H16").Value = 28.5
$wsSS.Range("I16").Value = 29.6
$wsSS.Range("D17").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("D18").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("E18").Value = 36
$wsSS.Range("F18").Value = 949
$wsSS.Range("G18").Value = 1226
$wsSS.Range("I18").Value = 34.1
$wsSS.Range("D19").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("E19").Value = 36
$wsSS.Range("F19").Value = 1190
$wsSS.Range("G19").Value = 1090
$wsSS.Range("H19").Value = 33.1
$wsSS.Range("I19").Value = 30.3
$wsSS.Range("D20").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("D21").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("E21").Value = 36
$wsSS.Range("F21").Value = 1274
$wsSS.Range("G21").Value = 1086
$wsSS.Range("H21").Value = 35.4
$wsSS.Range("I21").Value = 30.2
$wsSS.Range("D22").Value = '2025-12-20T17:00:00Z'
$wsSS.Range("D23").Value = '2025-12-20T17:00:00Z'

# ---- Meta_ext ----
$wsMeta.Range("B2").Value = '2025-12-20T17:00:00Z'
$wsMeta.Range("D2").Value = 75

